# PRM6: norm(similarity) + position
# Fills in a new experiment row (row 9) describing a "prm5_resnet50" run
# and adds a remark to the previous run (row 8, column J), plus the
# bookkeeping UI state (selected cell, column width) that Excel recorded
# when the author made the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: new PRM6 run (batchsize 4 x 16, l1norm similarity, Y/Y
#     normalize/position, 16 machines, temp4, 64 FP, prm5_resnet50) ---

# A9 gets the same "RUNNING" status styling (colored font) already used in
# A6/A7/A8; copy that formatting first, then set the text.
$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A9").Value = "RUNNING"

$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = "l1norm"
$ws.Range("E9").Value = "Y"
$ws.Range("F9").Value = "Y"
$ws.Range("G9").Value = 16
$ws.Range("K9").Value = "prm5_resnet50"
$ws.Range("H9").Value = "temp4"
$ws.Range("I9").Value = 64
# J9 stays blank, matching the rest of the table.

# --- Row 8: add a remark in column J (red warning text, like J5) ---
$ws.Range("J8").Value = "too lower; conv slow"
$ws.Range("J8").Font.Color = 255

# --- UI state captured by Excel on save ---
$ws.Range("D13").Select()
$ws.Columns(10).ColumnWidth = 17.666666666666668
